$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.336.12"
$ws.Range("E2").Value = "  -3.24%  "

$ws.Range("D3").Value = "2.363.24"
$ws.Range("E3").Value = "  -4.66%  "

$ws.Range("D5").Value = "'310.41"
$ws.Range("E5").Value = "  -2.83%  "

$ws.Range("D6").Value = "'85.22"
$ws.Range("E6").Value = "  -7.99%  "

$ws.Range("E7").Value = "  -3.47%  "

$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").Value = "'0.488"
$ws.Range("E9").Value = "  -4.77%  "

$ws.Range("D10").Value = "'0.0826"
$ws.Range("E10").Value = "  -4.56%  "

$ws.Range("D11").Value = "'30.14"
$ws.Range("E11").Value = "  -9.14%  "

$ws.Range("D12").Value = "'0.109"
$ws.Range("E12").Value = "  -0.91%  "

$ws.Range("D13").Value = "2.735.39"
$ws.Range("E13").Value = "  -4.40%  "

$ws.Range("D14").Value = "'6.45"
$ws.Range("E14").Value = "  -6.48%  "

$ws.Range("D15").Value = "'14.84"
$ws.Range("E15").Value = "  -4.71%  "

$ws.Range("D16").Value = "2.358.54"
$ws.Range("E16").Value = "  -4.38%  "

$ws.Range("D17").Value = "'0.753"
$ws.Range("E17").Value = "  -5.33%  "

$ws.Range("D18").Value = "40.379.67"
$ws.Range("E18").Value = "  -3.02%  "

$ws.Range("D19").Value = "0.0₃0905"
$ws.Range("E19").Value = "  -4.02%  "

$ws.Range("D20").Value = "'6.10"
$ws.Range("E20").Value = "  -5.53%  "

$ws.Range("D21").Value = "'68.02"
$ws.Range("E21").Value = "  -4.04%  "

$ws.Range("D22").Value = "'10.73"
$ws.Range("E22").Value = "  -4.62%  "

$ws.Range("D23").Value = "'234.61"
$ws.Range("E23").Value = "  -2.15%  "

$ws.Range("E24").Value = "  -6.38%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").Value = "'1.79"
$ws.Range("E26").Value = "  -8.06%  "

$ws.Range("D27").Value = "'23.57"
$ws.Range("E27").Value = "  -5.77%  "

$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  -1.12%  "

$ws.Range("D29").Value = "'9.19"
$ws.Range("E29").Value = "  -5.69%  "

$ws.Range("D30").Value = "'34.00"
$ws.Range("E30").Value = "  -7.29%  "

$ws.Range("D31").Value = "'152.79"
$ws.Range("E31").Value = "  -3.02%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").Value = "'5.17"
$ws.Range("E33").Value = "  -4.84%  "

$ws.Range("D34").Value = "'0.0726"
$ws.Range("E34").Value = "  -5.05%  "

$ws.Range("E35").Value = "  -5.46%  "

$ws.Range("E36").Value = "  -2.35%  "

$ws.Range("D37").Value = "'2.77"

$ws.Range("D38").Value = "'15.79"
$ws.Range("E38").Value = "  -8.07%  "

$ws.Range("D39").Value = "'0.0989"
$ws.Range("E39").Value = "  -4.77%  "

$ws.Range("D40").Value = "'1.69"
$ws.Range("E40").Value = "  -8.21%  "

$ws.Range("D41").Value = "'3.79"
$ws.Range("E41").Value = "  -5.51%  "

$ws.Range("D42").Value = "'2.38"
$ws.Range("E42").Value = "  -3.82%  "

$ws.Range("D43").Value = "1.962.67"
$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("D44").Value = "'0.0266"
$ws.Range("E44").Value = "  -6.46%  "

$ws.Range("D45").Value = "'17.61"
$ws.Range("E45").Value = "  -5.63%  "

$ws.Range("D46").Value = "'9.24"
$ws.Range("E46").Value = "  -2.95%  "

$ws.Range("D47").Value = "'2.68"
$ws.Range("E47").Value = "  -10.09%  "

$ws.Range("D48").Value = "2.602.49"
$ws.Range("E48").Value = "  -4.29%  "

$ws.Range("D49").Value = "'92.70"
$ws.Range("E49").Value = "  -5.40%  "

$ws.Range("D50").Value = "'71.10"
$ws.Range("E50").Value = "  -6.02%  "

$ws.Range("D51").Value = "'50.02"
$ws.Range("E51").Value = "  -4.40%  "
